$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '25.815.96'
$ws.Range('E2').Value = '  -1.31%  '
$ws.Range('D3').Value = '1.633.47'
$ws.Range('E3').Value = '  -1.39%  '
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '1.002'
$ws.Range('E4').Value = '  -0.39%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '214.63'
$ws.Range('E5').Value = '  -0.41%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '0.5024'
$ws.Range('E6').Value = '  -1.69%  '
$ws.Range('E7').Value = '  -0.39%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.2566'
$ws.Range('E8').Value = '  -0.62%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.06384'
$ws.Range('E9').Value = '  -0.54%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '19.70'
$ws.Range('E10').Value = '  -1.19%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.07678'
$ws.Range('E11').Value = '  -1.58%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '4.251'
$ws.Range('E12').Value = '  -0.80%  '
$ws.Range('D13').Value = '1.634.65'
$ws.Range('E13').Value = '  -1.33%  '
$ws.Range('D14').Value = '1.858.48'
$ws.Range('E14').Value = '  -1.34%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '0.5428'
$ws.Range('E15').Value = '  -1.63%  '
$ws.Range('D16').Value = '0.0₅7933'
$ws.Range('E16').Value = '  -1.15%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '63.58'
$ws.Range('E17').Value = '  -0.77%  '
$ws.Range('D18').Value = '25.834.87'
$ws.Range('E18').Value = '  -1.27%  '
$ws.Range('E19').Value = '  -0.41%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '201.53'
$ws.Range('E20').Value = '  -3.77%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '4.338'
$ws.Range('E21').Value = '  -1.72%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '9.923'
$ws.Range('E22').Value = '  -1.43%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '5.959'
$ws.Range('E23').Value = '  -0.97%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '1.003'
$ws.Range('E24').Value = '  -0.31%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '1.933'
$ws.Range('E25').Value = '  +11.25%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '141.93'
$ws.Range('E26').Value = '  -1.14%  '
$ws.Range('E27').Value = '  -3.21%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '15.66'
$ws.Range('E28').Value = '  -0.97%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '6.710'
$ws.Range('E29').Value = '  -3.87%  '
$ws.Range('E30').Value = '  -0.26%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '0.05000'
$ws.Range('E31').Value = '  -2.41%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '3.263'
$ws.Range('E32').Value = '  -2.69%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '3.187'
$ws.Range('E33').Value = '  -1.06%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '1.542'
$ws.Range('E34').Value = '  -1.50%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '2.367'
$ws.Range('E35').Value = '  -0.17%  '
$ws.Range('D36').Value = '1.171.59'
$ws.Range('E36').Value = '  +1.02%  '
$ws.Range('B37').Value = 'MXToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '2.634'
$ws.Range('E37').Value = '  -4.16%  '
$ws.Range('B38').Value = 'ARBITRUM'
$ws.Range('C38').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.8925'
$ws.Range('E38').Value = '  -3.89%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.5598'
$ws.Range('E39').Value = '  -1.67%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.01559'
$ws.Range('E40').Value = '  -1.93%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '1.002'
$ws.Range('E41').Value = '  -0.37%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '5.699'
$ws.Range('E42').Value = '  +1.14%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.8098'
$ws.Range('E43').Value = '  -2.93%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '99.65'
$ws.Range('E44').Value = '  -0.81%  '
$ws.Range('D45').Value = '1.770.41'
$ws.Range('E45').Value = '  -1.26%  '
$ws.Range('E46').Value = '  -1.65%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.4514'
$ws.Range('E47').Value = '  -0.72%  '
$ws.Range('E48').Value = '  -0.13%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '54.88'
$ws.Range('E49').Value = '  -1.56%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.05077'
$ws.Range('E50').Value = '  +0.68%  '
$ws.Range('B51').Value = 'USDD'
$ws.Range('C51').Value = 'https://coinranking.com/coin/z2PZIKQL7+usdd-usdd'
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '1.001'
$ws.Range('E51').Value = '  -0.46%  '
